$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the s1cDNADate value: "1.10.20" -> "01.10.20" for every row (D2:D27)
#    A literal text assignment here gets re-interpreted as a date (since the
#    text looks like D.M.YY), which would silently swap in a date number
#    format and a new style. Route it through a text formula, then paste
#    the computed value back as a value-only paste: that keeps the result a
#    plain shared-string literal without ever touching NumberFormat/style.
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $cell.Formula = "=""01.10.20"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# 2) Replace the "No" literal in column I (roboticS1Prep) with a FALSE()
#    formula, formatted to display as TRUE/FALSE, for every data row.
for ($r = 2; $r -le 27; $r++) {
    $cell = $ws.Cells.Item($r, 9)
    $cell.NumberFormat = '"TRUE";"TRUE";"FALSE"'
    $cell.Formula = "=FALSE()"
}

# 3) Move the active selection from H2:H27 to D27.
$ws.Range("D27").Select()
